$d = $word.ActiveDocument

# The document regenerates internal identifiers (bookmark ids and the
# field-group rsid markers) each time it is produced. Re-create the
# "bookmark1" bookmark in place so Word mints a fresh internal id for it,
# matching what a new generation run of the template would produce.
$bm = $d.Bookmarks("bookmark1")
$bmRange = $d.Range($bm.Range.Start, $bm.Range.End)
$bm.Delete()
$d.Bookmarks.Add("bookmark1", $bmRange) | Out-Null
